$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.13"
$ws.Range("E2").Value = "'0.72%"
$ws.Range("D3").Value = "'29.41"
$ws.Range("E3").Value = "'8.19%"
$ws.Range("D4").Value = "'5.186"
$ws.Range("E4").Value = "'1.55%"
$ws.Range("D5").Value = "'0.05735"
$ws.Range("E5").Value = "'0.82%"
$ws.Range("D6").Value = "'6.550"
$ws.Range("E6").Value = "'0.37%"
$ws.Range("D7").Value = "'3.097"
$ws.Range("E7").Value = "'2.94%"
$ws.Range("D8").Value = "'0.8581"
$ws.Range("E8").Value = "'4.73%"
$ws.Range("D9").Value = "'0.8680"
$ws.Range("E9").Value = "'0.89%"
$ws.Range("D10").Value = "'0.1364"
$ws.Range("E10").Value = "'2.23%"
$ws.Range("D11").Value = "'0.07086"
$ws.Range("E11").Value = "'2.00%"
$ws.Range("D12").Value = "'0.02980"
$ws.Range("E12").Value = "'4.43%"
$ws.Range("D13").Value = "'0.09385"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'0.77%"
$ws.Range("D15").Value = "'0.04135"
$ws.Range("E15").Value = "'1.57%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006022"
$ws.Range("E16").Value = "'-3.09%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007489"
$ws.Range("E17").Value = "'5,082.83%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.490"
$ws.Range("E18").Value = "'-0.48%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.274"
$ws.Range("E19").Value = "'-1.83%"
$ws.Range("B20").Value = "One"
$ws.Range("C20").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D20").Value = "'0.01028"
$ws.Range("E20").Value = "'1,615.87%"
$ws.Range("D22").Value = "'0.03397"
$ws.Range("E22").Value = "'5.42%"
$ws.Range("D23").Value = "'0.1303"
$ws.Range("E23").Value = "'2.29%"
$ws.Range("D24").Value = "'3.469"
$ws.Range("E24").Value = "'-2.34%"
$ws.Range("D25").Value = "'0.1380"
$ws.Range("E25").Value = "'0.44%"
$ws.Range("D26").Value = "'0.005008"
$ws.Range("E26").Value = "'12.09%"
$ws.Range("D27").Value = "'0.001225"
$ws.Range("E27").Value = "'0.84%"
$ws.Range("E28").Value = "'2.55%"
$ws.Range("D40").Value = "'0.03752"
$ws.Range("E40").Value = "'0.77%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005719"
$ws.Range("E41").Value = "'-4.34%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E42").Value = "'1.34%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002428"
$ws.Range("E43").Value = "'5.58%"
$ws.Range("D44").Value = "'0.008487"
$ws.Range("E44").Value = "'-12.63%"
$ws.Range("D45").Value = "'0.00005259"
$ws.Range("E45").Value = "'3.14%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D47").Value = "'0.05699"
$ws.Range("E47").Value = "'-43.56%"
$ws.Range("E48").Value = "'-9.29%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.00%"
